# Updated symbol list on Mon Dec 26 11:47:27 UTC 2022 with GitHub Actions
# Refresh the scraped coin prices (column D) and a couple of "Best/Worst in 24h"
# suffix labels (column E) to the latest pulled values.
#
# Column D cells store price text that looks numeric (e.g. "243.53", "5.410",
# "0.00000000751") but must stay as literal text (matches the workbook's
# original inlineStr/string cells, preserving exact digits/trailing zeros).
# Prefixing with a leading apostrophe forces Excel to keep the literal text
# instead of auto-converting to a number; resetting the Style back to
# "Normal" afterwards clears the "quote prefix" cell style Excel applies,
# so the cell's style index is left unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'243.53"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Value = "'23.04"
$ws.Range("D3").Style = "Normal"
$ws.Range("D4").Value = "'5.410"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Value = "'0.05917"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Value = "'3.457"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Value = "'6.554"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Value = "'0.8119"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Value = "'0.9133"
$ws.Range("D9").Style = "Normal"
$ws.Range("D11").Value = "'0.07422"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Value = "'0.03264"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Value = "'0.03067"
$ws.Range("D13").Style = "Normal"
$ws.Range("D16").Value = "'0.001571"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Value = "'0.04678"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Value = "'0.0005944"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Value = "'0.005943"
$ws.Range("D19").Style = "Normal"
$ws.Range("D21").Value = "'0.0009851"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Value = "'0.00008605"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Value = "'3.617"
$ws.Range("D23").Style = "Normal"
$ws.Range("D25").Value = "'0.3241"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Value = "'0.1325"
$ws.Range("D26").Style = "Normal"
$ws.Range("D40").Value = "'0.03963"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Value = "'0.006198"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "40KickTokenKICK"
$ws.Range("D42").Value = "'0.1077"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Value = "'0.002552"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Value = "'0.008617"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "43LocalTradersLCTBestin24h"
$ws.Range("D45").Value = "'0.00005172"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Value = "'0.00000000751"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Value = "'0.8796"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Value = "'0.002265"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Value = "'0.00002102"
$ws.Range("D49").Style = "Normal"
